# 2015_MMMILC_dataCleaningLog.xlsx - "updated data cleaning log"
# Append one new data-cleaning-log entry (obs.ID 7945) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = 7945
$ws.Range("B53").Value = "stage.length 5mm"
$ws.Range("C53").Value = "changed to L1-5, because L1 most likely"

# Move the selection to the newly added row, matching where the author
# was last working in the sheet.
$ws.Range("C53").Select()
